# ----------------------------------------------------------------------------
# Boletin Epi Pereira - poisson.xlsx - semana 26 de 2025
#
# The source surveillance export inserted a new event (420 - Leishmaniasis
# cutanea) ahead of 455 - Leptospirosis, which pushed every subsequent event
# code/name down by one row (through 591 - Vigilancia integrada ...). This
# block re-keys A20:B27 to the new alignment, then every row's Esperado/
# Observado/valor p (columns C/D/E) is refreshed with this week's figures.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-key the event code/name column for the inserted row (420) and the
#     cascading shift of every event below it (rows 20-27). ---
$eventRelabel = @(
    @{Row=20; Evento="420"; Nombre="Leishmaniasis cutanea"},
    @{Row=21; Evento="455"; Nombre="Leptospirosis"},
    @{Row=22; Evento="465"; Nombre="Malaria"},
    @{Row=23; Evento="535"; Nombre="Meningitis bacteriana y enfermedad meningoc”cica"},
    @{Row=24; Evento="549"; Nombre="Morbilidad materna extrema"},
    @{Row=25; Evento="560"; Nombre="Mortalidad perinatal y neonatal tardia"},
    @{Row=26; Evento="580"; Nombre="Mortalidad por dengue"},
    @{Row=27; Evento="591"; Nombre="Vigilancia integrada de muertes en menores de cinco anos por infeccion respiratoria aguda - enfermedad diarreica aguda y/o desnutricion"}
)
foreach ($item in $eventRelabel) {
    $cell = $ws.Cells.Item($item.Row, 1)
    $cell.NumberFormat = "@"   # keep the numeric-looking code stored as text
    $cell.Value = $item.Evento
    $ws.Cells.Item($item.Row, 2).Value = $item.Nombre
}

# --- Refresh Esperado (C) / Observado (D) / valor p (E) for semana 26. ---
$statUpdates = @(
    @{Row=2; D=0; E=1},
    @{Row=3; C=1; D=6; E=0},
    @{Row=4; C=0; D=0; E=1},
    @{Row=5; C=6; D=10; E=0.04},
    @{Row=6; C=2; D=3; E=0.18},
    @{Row=7; C=1; D=5; E=0},
    @{Row=9; C=48},
    @{Row=10; C=1; E=0.37},
    @{Row=11; D=4; E=0.02},
    @{Row=12; D=5; E=0.1},
    @{Row=13; C=42},
    @{Row=14; C=4; E=0.02},
    @{Row=15; C=2; E=0.14},
    @{Row=16; D=0; E=1},
    @{Row=17; D=8; E=0.14},
    @{Row=18; C=1; E=0.37},
    @{Row=19; C=6; D=0; E=0},
    @{Row=20; C=0; D=0; E=1},
    @{Row=21; C=2; D=0; E=0.14},
    @{Row=23; C=0; D=0; E=1},
    @{Row=24; C=6; D=5; E=0.16},
    @{Row=25; C=1; D=1; E=0.37},
    @{Row=27; C=0; E=1},
    @{Row=30; C=2; E=0.27},
    @{Row=31; D=0; E=1},
    @{Row=32; D=2; E=0.02},
    @{Row=33; C=9; E=0},
    @{Row=34; C=7; D=7; E=0.15}
)
foreach ($item in $statUpdates) {
    if ($item.ContainsKey("C")) { $ws.Cells.Item($item.Row, 3).Value = $item.C }
    if ($item.ContainsKey("D")) { $ws.Cells.Item($item.Row, 4).Value = $item.D }
    if ($item.ContainsKey("E")) { $ws.Cells.Item($item.Row, 5).Value = $item.E }
}
